# "fix submission set author order"
# The RefID column (A) is renumbered and the Authors column (C) is
# populated with the correctly-ordered author list for each publication
# row. The underlying row order (by Title/Journal/Year/Accession/PMID)
# does not change - only the RefID values and the Authors text change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- RefID (A) + Authors (C) updates, row by row ---

$ws.Range("A2").Value = 6

$ws.Range("A3").Value = 8

$ws.Range("A4").Value = 11
$ws.Range("C4").Value = 'Azuero O., Ou T., Lefrancq N., Nikolay B., Mckee C., Cappelle J., Hul V., Hoem T., Lemey P., Rahman M., Islam A., Gurley E., Hul V., Hoem T., Heng O., Williams D., Cappelle J., Salje H., Duong V.'

$ws.Range("A5").Value = 13
$ws.Range("C5").Value = 'Chang L., Rahman S., Hassan S., Olival K., Mohamed M., Hassan L., Saad N., Shohaimi S., Mamat Z., Naim M., Epstein J., Suri A., Field H., Daszak P., Henipavirus ecology research G.'

$ws.Range("A6").Value = 15
$ws.Range("C6").Value = 'Duong V., Cappelle J., Hul V., Buchy P.'

$ws.Range("A7").Value = 20
$ws.Range("C7").Value = 'Honko A., Johnson J., Hensley L., Wohl S., Barnes K., Sabeti P., Olinger G., Jahrling P., Wohl S., Barnes K., Sabeti P., Hensley L., Olinger G., Jahrling P., Sword J., Honko A.'

$ws.Range("A8").Value = 21
$ws.Range("C8").Value = 'Kohl C., Siriwardana S., Muzeniek T., Perera T., Bas D., Oeruc M., Brinkmann A., Becker-ziaja B., Schwarz F., Jeevatharan H., Weerasena J., Handunnetti S., Perera I., Premawansa G., Premawansa S., Yapa W., Nitsche A.'

$ws.Range("A9").Value = 24
$ws.Range("C9").Value = 'Sendow I., Ratnawati A., Taylor T., Abdul adjid R., Saepulloh M., Barr J., Daniels P., Field H.'

$ws.Range("A10").Value = 25
$ws.Range("C10").Value = 'Sharifah S., Sohayati A., Maizan M., Chang L., Sharina M., Syamsiah A., Latiffah K., Siti suri A., Zaini C., Humes F., Daszak P., Epstein J.'

$ws.Range("A11").Value = 29
$ws.Range("C11").Value = 'Wacharapluesadee S., Ngamprasertwong T., Supavonwong P., Phumesin P., Ratanasetyuth N., Boongird K., Wanghongsa S., Supavonwong P., Phumesin P., Ratanasetyuth N., Boongird K., Wanghongsa S., Hemachudha T.'

$ws.Range("A12").Value = 30
$ws.Range("C12").Value = 'Whitmer S., Lo M., Sazzad H., Zufan S., Gurley E., Sultana S., Amman B., Ladner J., Rahman M., Doan S., Satter S., Flora M., Montgomery J., Nichol S., Spiropoulou C., Lo M., Zufan S., Nichol S., Spiropoulou C., Klena J.'

$ws.Range("A13").Value = 36
$ws.Range("C13").Value = 'Rahman M., Miah M., Hossain M., Satter S., Klena J., Shirin T., Montgomery J., Rahman D., Rahman S., Miah M., Rahman S., Rahman D., Hossain M., Satter S., Klena J., Shirin T., Montgomery J., Rahman M.'

$ws.Range("A14").Value = 38
$ws.Range("C14").Value = 'Wacharapluesadee S., Lumlertdacha B., Boongird K., Wanghongsa S., Chanhome L., Rollin P., Stockton P., Rupprecht C., Ksiazek T., Hemachudha T.'

$ws.Range("A15").Value = 40

# --- Re-apply the sort on the RefID column so the sheet records a
#     sortState (the data is already in ascending RefID order, so this
#     does not move any rows - it only records that the range was last
#     sorted by column A). ---

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A15"))
$ws.Sort.SetRange($ws.Range("A1:G15"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- Final selection left on the Authors column ---

$ws.Range("C2:C15").Select()
